$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GIS")

# Row 4: Inventory
$ws.Range("B4").Value = 1759000000.0
$ws.Range("C4").Value = 1713000000.0
$ws.Range("D4").Value = 1605000000.0
$ws.Range("E4").Value = 1426000000.0
$ws.Range("F4").Value = 1542000000.0

# Row 14: Accounts Payable
$ws.Range("B14").Value = 3392000000.0
$ws.Range("C14").Value = 3399000000.0
$ws.Range("D14").Value = 3184000000.0
$ws.Range("E14").Value = 3248000000.0
$ws.Range("F14").Value = 2932000000.0

# Row 19: Long Term Tax Liability (Deferred)
$ws.Range("B19").Value = 2006000000.0
$ws.Range("C19").Value = 1939000000.0
$ws.Range("D19").Value = 1925000000.0
$ws.Range("E19").Value = 1947000000.0
$ws.Range("F19").Value = 2027000000.0
